# Refit the Beta Abundance model and refresh the parameter-estimate
# panels (posterior mean and the Pr(x>0)/Pr(x<0) tables) that feed the
# RSquared calculation used when plotting this panel.
$wb = $excel.ActiveWorkbook

$dataMean = @(
    @(-6.2659309253239, -0.5988942842363025, -0.0140986211690393, -0.395317187318695, 0.007248249644022577, 0.1951207373002259, 0.01087754231121822, -0.07928421786462415),
    @(-10.72164466035066, -0.8043154816306575, 0.04124935672767511, -0.8350964202807636, -0.03432429192355456, 0.2648233406737049, 0.03979806246669471, -0.2832089387179023),
    @(2.043122287843016, -0.1112767069613757, -0.05932203237308012, 0.06625336085183849, -0.07003508609985085, -0.338410543024619, -0.06037784303232136, 1.65365699251615),
    @(-9.508006568357688, -0.1681438044433469, -0.01522853408290514, -0.6748524692989257, -0.07290528571770309, 0.04581744617959855, 0.05543605725446554, 0.1340567878675038),
    @(-10.27501837697094, 0.4121991287101447, -0.01259309166135484, -0.01814947604667764, 0.0259453283256228, 0.2227597874779113, 0.02026743534997227, 0.2977731179121011),
    @(-6.398771443458977, 1.195711469806844, -0.0004365683546076519, 0.1833756311464365, 0.007098102094892961, 0.07520610624885114, 0.01843277933350791, 0.110225409380142),
    @(-8.414319109939104, 1.713710700555608, -0.001685712123157412, 0.3518928857056302, 0.005271989081050034, 0.06802620777648208, 0.02866923596795208, 0.324948894735677),
    @(-5.003515765543416, 2.094526694466488, -0.01313730066515563, 0.685502183251541, 0.003223352991792302, 0.07947413272493827, 0.00889948454518103, 0.1126379670019341),
    @(-5.269847661638758, 0.9520113185559916, 0.00003830684377228642, 0.2634160282476582, 0.005670951710796803, -0.04569753286671877, 0.01000221892351058, 0.1419500704792796),
    @(-0.08081600901010699, 0.1589588998973643, -0.01393218973184394, 0.3801712135025674, -0.001320550156367274, 0.08587700128372658, -0.01073643459022128, -0.5385430273393251),
    @(-3.737562131959649, 1.198938645820635, -0.005669955498870197, 0.4994194862424781, -0.002439953630670999, -0.02565983010042778, -0.000219298539083331, 0.0526718488413991),
    @(-6.818080241163051, -1.542599298827969, 0.01911266511003648, -0.2324557779673498, 0.024870453764421, 0.0434807513647179, -0.009920020765254346, 0.485676954731076)
)
$ws = $wb.Worksheets.Item("Posterior mean")
for ($r = 0; $r -lt $dataMean.Length; $r++) {
    $row = $dataMean[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item(2 + $r, 2 + $c).Value2 = $row[$c]
    }
}

$dataPosProb = @(
    @(0.04284375, 0.08712499999999999, 0.0043125, 0.05775, 0.919875, 0.9770624999999999, 0.7776875, 0.28696875),
    @(0.008718750000000001, 0.0068125, 1, 0.07178125, 0.000125, 0.9799687500000001, 0.99290625, 0.009124999999999999),
    @(0.69171875, 0.36009375, 0, 0.52209375, 0, 0, 0, 1),
    @(0.00009375, 0.24621875, 0, 0.19190625, 0, 0.79925, 1, 0.87790625),
    @(0.00009375, 0.84403125, 0.00128125, 0.4939375, 1, 0.9975000000000001, 0.9672500000000001, 0.99865625),
    @(0.051, 0.996625, 0.457375, 0.78303125, 0.9649375, 0.7708125, 0.90790625, 0.80221875),
    @(0.01546875, 1, 0.304625, 0.94271875, 0.9545625, 0.7409375, 0.97975, 0.999375),
    @(0.1104375, 1, 0.000125, 0.985375, 0.81459375, 0.7610625, 0.735625, 0.79571875),
    @(0.140375, 0.90821875, 0.502, 0.75125, 0.7905624999999999, 0.384125, 0.7185, 0.737125),
    @(0.473, 0.57984375, 0.10484375, 0.7761875, 0.46271875, 0.6921875, 0.3445, 0.01484375),
    @(0.2384375, 0.9348125, 0.3039375, 0.839, 0.42525, 0.4446875, 0.50109375, 0.5748124999999999),
    @(0.11515625, 0.03646875, 0.9975937499999999, 0.31625, 1, 0.615125, 0.319625, 0.99684375)
)
$ws = $wb.Worksheets.Item("Pr(x>0)")
for ($r = 0; $r -lt $dataPosProb.Length; $r++) {
    $row = $dataPosProb[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item(2 + $r, 2 + $c).Value2 = $row[$c]
    }
}

$dataNegProb = @(
    @(0.95715625, 0.912875, 0.9956874999999999, 0.94225, 0.080125, 0.0229375, 0.2223125, 0.71303125),
    @(0.99128125, 0.9931875, 0, 0.92821875, 0.999875, 0.02003125, 0.00709375, 0.990875),
    @(0.30828125, 0.63990625, 1, 0.47790625, 1, 1, 1, 0),
    @(0.99990625, 0.75378125, 1, 0.8080937500000001, 1, 0.20075, 0, 0.12209375),
    @(0.99990625, 0.15596875, 0.99871875, 0.5060625, 0, 0.0025, 0.03275, 0.00134375),
    @(0.949, 0.003375, 0.542625, 0.21696875, 0.0350625, 0.2291875, 0.09209375, 0.19778125),
    @(0.98453125, 0, 0.695375, 0.05728125, 0.0454375, 0.2590625, 0.02025, 0.000625),
    @(0.8895625, 0, 0.999875, 0.014625, 0.18540625, 0.2389375, 0.264375, 0.20428125),
    @(0.859625, 0.09178124999999999, 0.498, 0.24875, 0.2094375, 0.615875, 0.2815, 0.262875),
    @(0.527, 0.42015625, 0.89515625, 0.2238125, 0.53728125, 0.3078125, 0.6555, 0.98515625),
    @(0.7615625, 0.0651875, 0.6960625, 0.161, 0.57475, 0.5553125, 0.49890625, 0.4251875),
    @(0.88484375, 0.96353125, 0.00240625, 0.68375, 0, 0.384875, 0.680375, 0.00315625)
)
$ws = $wb.Worksheets.Item("Pr(x<0)")
for ($r = 0; $r -lt $dataNegProb.Length; $r++) {
    $row = $dataNegProb[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item(2 + $r, 2 + $c).Value2 = $row[$c]
    }
}
